$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.456.02"
$ws.Range("E2").Value = "  +4.84%  "

$ws.Range("D3").Value = "3.355.34"
$ws.Range("E3").Value = "  +9.55%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'255.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.41%  "

$ws.Range("D6").Value = "'623.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.06%  "

$ws.Range("D7").Value = "'1.18"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.72%  "

$ws.Range("D8").Value = "'0.385"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.02%  "

$ws.Range("D9").Value = "'0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").Value = "3.350.66"
$ws.Range("E10").Value = "  +9.68%  "

$ws.Range("D11").Value = "'0.802"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.21%  "

$ws.Range("D12").Value = "'0.199"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.60%  "

$ws.Range("D13").Value = "98.106.35"
$ws.Range("E13").Value = "  +5.02%  "

$ws.Range("D14").Value = "'35.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.44%  "

$ws.Range("E15").Value = "  +3.30%  "

$ws.Range("D16").Value = "3.972.98"
$ws.Range("E16").Value = "  +9.77%  "

$ws.Range("E17").Value = "  +4.68%  "

$ws.Range("D18").Value = "3.356.22"
$ws.Range("E18").Value = "  +10.01%  "

$ws.Range("E19").Value = "  +2.52%  "

$ws.Range("D20").Value = "'14.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.10%  "

$ws.Range("D21").Value = "'485.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +11.53%  "

$ws.Range("E22").Value = "  +3.80%  "

$ws.Range("D23").Value = "'0.0000206"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +10.39%  "

$ws.Range("D24").Value = "'9.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.91%  "

$ws.Range("D25").Value = "'5.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.05%  "

$ws.Range("D26").Value = "'88.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.63%  "

$ws.Range("D27").Value = "'11.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.47%  "

$ws.Range("D28").Value = "3.531.38"
$ws.Range("E28").Value = "  +10.47%  "

$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("B30").Value = "Cronos"
$ws.Range("C30").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D30").Value = "'0.185"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.98%  "

$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.248"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.89%  "

$ws.Range("D32").Value = "'0.125"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.84%  "

$ws.Range("D33").Value = "'0.996"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -11.04%  "

$ws.Range("D34").Value = "'9.26"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.55%  "

$ws.Range("D35").Value = "'27.27"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.62%  "

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D36").Value = "'7.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.85%  "

$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.151"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.19%  "

$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "'515.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.34%  "

$ws.Range("E39").Value = "  +3.30%  "

$ws.Range("E40").Value = "  +3.82%  "

$ws.Range("E41").Value = "  +4.08%  "

$ws.Range("E42").Value = "  +2.25%  "

$ws.Range("D43").Value = "'3.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.18%  "

$ws.Range("D44").Value = "'3.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.52%  "

$ws.Range("D46").Value = "'0.776"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +17.90%  "

$ws.Range("D47").Value = "'161.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.90%  "

$ws.Range("D48").Value = "'1.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.71%  "

$ws.Range("D49").Value = "'1.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.81%  "

$ws.Range("D50").Value = "'45.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.31%  "

$ws.Range("D51").Value = "'4.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.63%  "
